{"js": "// Update the worked \"two-digit \u00f7 one-digit\" answers in the practice\n// table. The sheet has one table; only every 4th row (0, 4, 8, 12, 16)\n// holds the five answer cells per line, the rows in between are blank\n// (left for handwriting). Replace each answer cell's text in place so\n// the existing run formatting (TimeNewRoman, sz 30) is preserved.\n\nconst table = context.document.body.tables.getFirst();\n\n// row index -> new cell values (left to right)\nconst newAnswers = {\n  0: [\"32\u00f76=5, 2\", \"71\u00f77=10, 1\", \"61\u00f78=7, 5\", \"25\u00f75=5, 0\", \"70\u00f72=35, 0\"],\n  4: [\"58\u00f74=14, 2\", \"49\u00f76=8, 1\", \"24\u00f77=3, 3\", \"92\u00f78=11, 4\", \"60\u00f76=10, 0\"],\n  8: [\"94\u00f75=18, 4\", \"57\u00f73=19, 0\", \"80\u00f72=40, 0\", \"60\u00f76=10, 0\", \"37\u00f75=7, 2\"],\n  12: [\"58\u00f74=14, 2\", \"20\u00f77=2, 6\", \"77\u00f73=25, 2\", \"30\u00f79=3, 3\", \"57\u00f78=7, 1\"],\n  16: [\"74\u00f73=24, 2\", \"31\u00f75=6, 1\", \"55\u00f74=13, 3\", \"62\u00f76=10, 2\", \"31\u00f73=10, 1\"],\n};\n\nfor (const rowIndex of Object.keys(newAnswers)) {\n  const values = newAnswers[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(Number(rowIndex), col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worked \"two-digit \u00f7 one-digit\" answers in the practice\n# table. The sheet has one table; only every 4th row (rows 1, 5, 9, 13,\n# 17 in 1-based COM indexing) holds the five answer cells per line, the\n# rows in between are blank (left for handwriting). Assigning directly\n# to Cell.Range.Text replaces just the run text and keeps the existing\n# character formatting (TimeNewRoman, sz 30) intact.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newAnswers = @{\n    1  = @(\"32\u00f76=5, 2\", \"71\u00f77=10, 1\", \"61\u00f78=7, 5\", \"25\u00f75=5, 0\", \"70\u00f72=35, 0\")\n    5  = @(\"58\u00f74=14, 2\", \"49\u00f76=8, 1\", \"24\u00f77=3, 3\", \"92\u00f78=11, 4\", \"60\u00f76=10, 0\")\n    9  = @(\"94\u00f75=18, 4\", \"57\u00f73=19, 0\", \"80\u00f72=40, 0\", \"60\u00f76=10, 0\", \"37\u00f75=7, 2\")\n    13 = @(\"58\u00f74=14, 2\", \"20\u00f77=2, 6\", \"77\u00f73=25, 2\", \"30\u00f79=3, 3\", \"57\u00f78=7, 1\")\n    17 = @(\"74\u00f73=24, 2\", \"31\u00f75=6, 1\", \"55\u00f74=13, 3\", \"62\u00f76=10, 2\", \"31\u00f73=10, 1\")\n}\n\nforeach ($rowIndex in $newAnswers.Keys) {\n    $values = $newAnswers[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
